$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "「彼を覚えていますか」" (row 196). All rows below it
# shift up by one, matching the author's diff (dimension A1:C285 -> A1:C284).
$ws.Rows.Item(196).Delete()
